# Updated cryptos list on Mon Apr 29 19:55:56 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'62.941.67"
$ws.Range("E2").Value = "  -1.41%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.178.79"
$ws.Range("E3").Value = "  -4.26%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'590.96"
$ws.Range("E5").Value = "  -2.33%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'134.36"
$ws.Range("E6").Value = "  -5.87%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "'3.178.04"
$ws.Range("E8").Value = "  -4.26%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.88%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.141"
$ws.Range("E10").Value = "  -6.44%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'5.23"
$ws.Range("E11").Value = "  -5.99%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -3.93%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "'0.0000236"
$ws.Range("E13").Value = "  -5.10%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'34.54"
$ws.Range("E14").Value = "  -1.50%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'3.695.35"
$ws.Range("E15").Value = "  -4.39%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  -1.19%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "'3.170.32"
$ws.Range("E17").Value = "  -4.46%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "'62.892.83"
$ws.Range("E18").Value = "  -1.61%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "'6.54"
$ws.Range("E19").Value = "  -4.78%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'460.16"
$ws.Range("E20").Value = "  -4.55%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'14.00"
$ws.Range("E21").Value = "  -0.85%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "'0.693"
$ws.Range("E22").Value = "  -6.29%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "'7.61"
$ws.Range("E23").Value = "  -4.74%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "'13.33"
$ws.Range("E24").Value = "  -4.64%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'82.66"
$ws.Range("E25").Value = "  -2.92%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.01%  "

# Row 27 - FirstDigitalUSD
$ws.Range("E27").Value = "  +0.06%  "

# Row 28 - PancakeSwap
$ws.Range("E28").Value = "  -4.13%  "

# Row 29 - RenderToken
$ws.Range("D29").Value = "'7.69"
$ws.Range("E29").Value = "  -6.88%  "

# Row 30 - NEARProtocol
$ws.Range("D30").Value = "'6.75"
$ws.Range("E30").Value = "  -5.93%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -6.25%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "'27.10"
$ws.Range("E32").Value = "  -6.38%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -4.02%  "

# Row 34 - Stacks
$ws.Range("D34").Value = "'2.36"
$ws.Range("E34").Value = "  -6.64%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  -6.81%  "

# Row 36 - Filecoin
$ws.Range("D36").Value = "'5.80"
$ws.Range("E36").Value = "  -4.89%  "

# Row 37 - OKB
$ws.Range("D37").Value = "'51.33"
$ws.Range("E37").Value = "  -2.16%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0705"
$ws.Range("E38").Value = "  -5.57%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "'0.0387"
$ws.Range("E39").Value = "  -3.23%  "

# Row 40 - Bittensor
$ws.Range("D40").Value = "'405.20"
$ws.Range("E40").Value = "  -6.84%  "

# Row 41 - Cosmos
$ws.Range("D41").Value = "'8.06"
$ws.Range("E41").Value = "  -3.76%  "

# Row 42 - Kaspa/dogwifhat swap
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.63"
$ws.Range("E42").Value = "  -4.78%  "

# Row 43 - dogwifhat/Kaspa swap
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.112"
$ws.Range("E43").Value = "  -5.38%  "

# Row 44 - Maker
$ws.Range("D44").Value = "'2.791.40"
$ws.Range("E44").Value = "  -10.65%  "

# Row 45 - TheGraph
$ws.Range("D45").Value = "'0.251"
$ws.Range("E45").Value = "  -6.54%  "

# Row 47 - Fetch.AI
$ws.Range("D47").Value = "'2.11"
$ws.Range("E47").Value = "  -6.20%  "

# Row 48 - Monero
$ws.Range("D48").Value = "'124.74"
$ws.Range("E48").Value = "  +0.34%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "'25.17"
$ws.Range("E49").Value = "  -4.65%  "

# Row 50 - Arweave
$ws.Range("D50").Value = "'34.39"
$ws.Range("E50").Value = "  -6.74%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -2.21%  "
